# Auto-generated Excel COM-interop script to apply the Pandaemonium_Profits.xlsx profit-recalculation update.
# For each touched cell: set the new numeric value, or ClearContents() when the cell should no longer exist.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 20954.4
$ws.Range("I125").Value = 1200
$ws.Range("J125").Value = 25893
$ws.Range("K125").Value = 10800
$ws.Range("L125").Value = 233037
$ws.Range("M125").Value = -8340
$ws.Range("N125").Value = -237957

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 993.17645
$ws.Range("I20").Value = 941.0714
$ws.Range("J20").Value = 1236.3334
$ws.Range("K20").Value = 941.0714
$ws.Range("L20").Value = 1236.3334
$ws.Range("M20").Value = -694.0714
$ws.Range("N20").Value = -1730.3334
$ws.Range("H35").Value = 36037
$ws.Range("J35").Value = 36037
$ws.Range("L35").Value = 36037
$ws.Range("N35").Value = -36657
$ws.Range("H86").Value = 1806.6364
$ws.Range("I86").Value = 1837.3
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1837.3
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -714.3
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 1806.6364
$ws.Range("I89").Value = 1837.3
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 9186.5
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -3570.5
$ws.Range("N89").Value = -18732
$ws.Range("H134").Value = 5817.7
$ws.Range("I134").Value = 5817.7
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 17453.1
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -14918.1
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 308.8889
$ws.Range("I22").Value = 274.6154
$ws.Range("K22").Value = 274.6154
$ws.Range("M22").Value = 75.38459999999998
$ws.Range("H31").Value = 3790.621
$ws.Range("I31").Value = 2069.2341
$ws.Range("J31").Value = 8048.7896
$ws.Range("K31").Value = 2069.2341
$ws.Range("L31").Value = 8048.7896
$ws.Range("M31").Value = -1774.2341
$ws.Range("N31").Value = -8638.7896
$ws.Range("H34").Value = 3790.621
$ws.Range("I34").Value = 2069.2341
$ws.Range("J34").Value = 8048.7896
$ws.Range("K34").Value = 2069.2341
$ws.Range("L34").Value = 8048.7896
$ws.Range("M34").Value = -1867.2341
$ws.Range("N34").Value = -8452.7896
$ws.Range("H38").Value = 5000
$ws.Range("J38").Value = 5000
$ws.Range("L38").Value = 5000
$ws.Range("N38").Value = -5754
$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5422
$ws.Range("H51").Value = 24839.6
$ws.Range("J51").Value = 24839.6
$ws.Range("L51").Value = 24839.6
$ws.Range("N51").Value = -26311.6
$ws.Range("H61").Value = 24839.6
$ws.Range("J61").Value = 24839.6
$ws.Range("L61").Value = 24839.6
$ws.Range("N61").Value = -25535.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8623275
$ws.Range("I5").Value = 493.23077
$ws.Range("K5").Value = 1479.69231
$ws.Range("M5").Value = -1367.69231
$ws.Range("H12").Value = 33333558
$ws.Range("I12").Value = 83333550
$ws.Range("J12").Value = 230.16667
$ws.Range("K12").Value = 250000650
$ws.Range("L12").Value = 690.50001
$ws.Range("M12").Value = -250000477
$ws.Range("N12").Value = -1036.50001
$ws.Range("H19").Value = 2833.6667
$ws.Range("I19").Value = 501
$ws.Range("K19").Value = 1503
$ws.Range("M19").Value = -1329
$ws.Range("H22").Value = 55556920
$ws.Range("I22").Value = 100000990
$ws.Range("J22").Value = 1826.25
$ws.Range("K22").Value = 300002970
$ws.Range("L22").Value = 5478.75
$ws.Range("M22").Value = -300002801
$ws.Range("N22").Value = -5816.75
$ws.Range("H27").Value = 55556920
$ws.Range("I27").Value = 100000990
$ws.Range("J27").Value = 1826.25
$ws.Range("K27").Value = 300002970
$ws.Range("L27").Value = 5478.75
$ws.Range("M27").Value = -300002868
$ws.Range("N27").Value = -5682.75
$ws.Range("H34").Value = 2739.2856
$ws.Range("J34").Value = 3042.8
$ws.Range("L34").Value = 9128.400000000001
$ws.Range("N34").Value = -9296.400000000001
$ws.Range("H105").Value = 6799.6665
$ws.Range("J105").Value = 6799.6665
$ws.Range("L105").Value = 20398.9995
$ws.Range("N105").Value = -25640.9995
$ws.Range("H112").Value = 2811.3572
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 2950.6924
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 8852.0772
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -11068.0772
$ws.Range("H131").Value = 513.9400000000001
$ws.Range("I131").Value = 278.09836
$ws.Range("J131").Value = 882.8205
$ws.Range("K131").Value = 834.2950800000001
$ws.Range("L131").Value = 2648.4615
$ws.Range("M131").Value = 4205.70492
$ws.Range("N131").Value = -12728.4615
$ws.Range("H135").Value = 8623275
$ws.Range("I135").Value = 493.23077
$ws.Range("K135").Value = 4439.07693
$ws.Range("M135").Value = -1904.07693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6830.615
$ws.Range("J80").Value = 3289.9
$ws.Range("L80").Value = 3289.9
$ws.Range("N80").Value = -5285.9
$ws.Range("H83").Value = 6830.615
$ws.Range("J83").Value = 3289.9
$ws.Range("L83").Value = 16449.5
$ws.Range("N83").Value = -26433.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H82").Value = 3456.2856
$ws.Range("I82").Value = 1673.5
$ws.Range("J82").Value = 5833.3335
$ws.Range("K82").Value = 1673.5
$ws.Range("L82").Value = 5833.3335
$ws.Range("M82").Value = -1312.5
$ws.Range("N82").Value = -6555.3335
$ws.Range("H85").Value = 3456.2856
$ws.Range("I85").Value = 1673.5
$ws.Range("J85").Value = 5833.3335
$ws.Range("K85").Value = 1673.5
$ws.Range("L85").Value = 5833.3335
$ws.Range("M85").Value = -425.5
$ws.Range("N85").Value = -8329.333500000001
$ws.Range("H132").Value = 3060.4285
$ws.Range("I132").Value = 2384.8
$ws.Range("J132").Value = 3961.2666
$ws.Range("K132").Value = 7154.400000000001
$ws.Range("L132").Value = 11883.7998
$ws.Range("M132").Value = -4624.400000000001
$ws.Range("N132").Value = -16943.7998
$ws.Range("H136").Value = 4152.5845
$ws.Range("I136").Value = 2596.738
$ws.Range("J136").Value = 6993.696
$ws.Range("K136").Value = 7790.214
$ws.Range("L136").Value = 20981.088
$ws.Range("M136").Value = -5240.214
$ws.Range("N136").Value = -26081.088

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8312.5
$ws.Range("H18").Value = 9500
$ws.Range("I18").Value = 9000
$ws.Range("K18").Value = 9000
$ws.Range("M18").Value = -8827
$ws.Range("H81").Value = 2299
$ws.Range("I81").Value = 2416.6667
$ws.Range("J81").Value = 2122.5
$ws.Range("K81").Value = 4833.3334
$ws.Range("L81").Value = 4245
$ws.Range("M81").Value = -3772.3334
$ws.Range("N81").Value = -6367
$ws.Range("H84").Value = 2299
$ws.Range("I84").Value = 2416.6667
$ws.Range("J84").Value = 2122.5
$ws.Range("K84").Value = 24166.667
$ws.Range("L84").Value = 21225
$ws.Range("M84").Value = -18862.667
$ws.Range("N84").Value = -31833
$ws.Range("H96").Value = 875
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 950
$ws.Range("K96").Value = 950
$ws.Range("L96").Value = 950
$ws.Range("M96").Value = 573
$ws.Range("N96").Value = -3696
